# Bỏ cột "Địa Chỉ Theo CCCD" (I) khỏi các dòng dữ liệu.
# Header (row 1) giữ nguyên, chỉ xóa nội dung của I2:I4.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2:I4").ClearContents()
